$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix tiny floating point differences on row 3 (recalculated values) ---
$ws.Range("C3").Value = 45406.92770666667
$ws.Range("D3").Value = 45406.92778109953

# --- Clear inherited column styling on the new plain (unstyled) cells ---
$plainCells = @("A4","B4","H4","A5","B5","H5","A6","B6","H6","I6")
foreach ($addr in $plainCells) {
    $ws.Range($addr).Style = "Normal"
}

# --- Copy the date/time and duration number formats from row 2 onto the new rows ---
$ws.Range("C2:D2").Copy()
$ws.Range("C4:D6").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("F2").Copy()
$ws.Range("F4:F6").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# --- Row 4 ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Snakes"
$ws.Range("C4").Value = 45406.93297440972
$ws.Range("D4").Value = 45406.93314144676
$ws.Range("F4").Value = 0.00016203703703703701
$ws.Range("H4").Value = "Team1"

# --- Row 5 ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Snakes"
$ws.Range("C5").Value = 45406.93342157407
$ws.Range("D5").Value = 45406.93350756945
$ws.Range("F5").Value = 0.00003472222222222222
$ws.Range("H5").Value = "Team3"

# --- Row 6 ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Snakes"
$ws.Range("C6").Value = 45406.94200854808
$ws.Range("D6").Value = 45406.94212322737
$ws.Range("F6").Value = 0.00010416666666666670
$ws.Range("H6").Value = "Team2"
$ws.Range("I6").Value = "Process10"
